# "settings set for Spasskaya Pad"
#
# 1. Insert a new column before AJ (column 36) for the new "mwindthrow"
#    parameter, shifting mdrought..biomasswoodfacb one column to the right.
# 2. Fill in the header + values for the new mwindthrow column.
# 3. Bump a couple of existing parameter values (H2, AF2).
# 4. Remove the stray formatting-only rows 59-64 below the data block.
# 5. Restore the intended column width for the freshly inserted column.
# 6. Update the active selection to match the edited view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new "mwindthrow" column at AJ ---------------------------
$ws.Columns("AJ:AJ").Insert()
$ws.Columns("AJ:AJ").ColumnWidth = 14.5

$ws.Range("AJ1").Value = "mwindthrow"
$ws.Range("AJ2").Value = 0.01
$ws.Range("AJ3").Value = 0.01
$ws.Range("AJ4").Value = 0.01
$ws.Range("AJ5").Value = 0.01
$ws.Range("AJ6").Value = 0.01

# --- tweak a couple of existing values -----------------------------------
$ws.Range("H2").Value = 16
$ws.Range("AF2").Value = 0.2

# --- drop the leftover formatting-only rows below the real data ----------
$ws.Rows("59:64").Delete()

# --- refresh the saved selection/view ------------------------------------
$ws.Range("AD8").Select()
